$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: set a cell to a literal *text* value (not auto-coerced to a number)
# while keeping/forcing the same "label" style already used elsewhere in the
# sheet for these placeholder cells (style index 14 -> general/text format).
# We do this by: forcing a text number format, assigning the text value, then
# pasting the cell-format only from an existing style-14 cell on top (this
# does not touch the value, only formatting/number format/style).
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text, $styleSourceAddr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Header strings (rich text runs inside shared strings) - update only the
# specific run's text via Characters(start, length) so the other runs/fonts
# in the same rich string stay untouched.
#   A8 = "Volume 30   Number  31"  -> "...Number  32"
#   C9 = "Report Covering the Week  7/31/2023  Through  8/6/2023"
#        -> "...8/7/2023  Through  8/13/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "32"

$ws.Range("C9").Characters(27, 9).Text = "8/7/2023"
$ws.Range("C9").Characters(47, 8).Text = "8/13/2023"

# ---------------------------------------------------------------------------
# Row 14 (Murder) - G14/H14 become blank-style placeholders ("0" / "***.*")
# ---------------------------------------------------------------------------
Set-TextCell $ws "G14" "0" "C14"
Set-TextCell $ws "H14" "***.*" "E14"

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
Set-TextCell $ws "D15" "0" "C14"
Set-TextCell $ws "E15" "***.*" "E14"
$ws.Range("L15").Value = -8.333333333333
$ws.Range("N15").Value = -68.571428571428

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 7
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = 96.078431372549
$ws.Range("M16").Value = -41.176470588235
$ws.Range("N16").Value = -83.552631578947

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 44
$ws.Range("H17").Value = 18.918918918918
$ws.Range("I17").Value = 287
$ws.Range("J17").Value = 249
$ws.Range("K17").Value = 15.261044176706
$ws.Range("L17").Value = 71.856287425149
$ws.Range("M17").Value = 51.052631578947
$ws.Range("N17").Value = -41.547861507128

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 90
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = 57.894736842105
$ws.Range("L18").Value = 119.512195121951
$ws.Range("M18").Value = -47.058823529411
$ws.Range("N18").Value = -91.847826086956

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 51.724137931034
$ws.Range("I19").Value = 251
$ws.Range("J19").Value = 216
$ws.Range("K19").Value = 16.203703703703
$ws.Range("L19").Value = 73.103448275862
$ws.Range("M19").Value = -1.953125
$ws.Range("N19").Value = -28.285714285714

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 166.666666666667
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 136.363636363636
$ws.Range("I20").Value = 93
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = 47.619047619047
$ws.Range("L20").Value = 121.428571428571
$ws.Range("M20").Value = -26.190476190476
$ws.Range("N20").Value = -88.848920863309

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 40.909090909090
$ws.Range("F21").Value = 133
$ws.Range("H21").Value = 40
$ws.Range("I21").Value = 842
$ws.Range("J21").Value = 656
$ws.Range("K21").Value = 28.353658536585
$ws.Range("L21").Value = 81.857451403887
$ws.Range("M21").Value = -10.330138445154
$ws.Range("N21").Value = -75.480489225393

# ---------------------------------------------------------------------------
# Row 23 (Transit)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C23" "0" "C14"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 60
$ws.Range("J23").Value = 36
$ws.Range("K23").Value = 66.666666666666
$ws.Range("L23").Value = 106.896551724138
$ws.Range("M23").Value = 130.769230769231

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 23
$ws.Range("E24").Value = -30.303030303030
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 5.217391304347
$ws.Range("I24").Value = 760
$ws.Range("J24").Value = 728
$ws.Range("K24").Value = 4.395604395604
$ws.Range("L24").Value = 55.419222903885
$ws.Range("M24").Value = -15.929203539823

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("E25").Value = 41.666666666666
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 43.181818181818
$ws.Range("I25").Value = 448
$ws.Range("J25").Value = 402
$ws.Range("K25").Value = 11.442786069651
$ws.Range("L25").Value = 40.438871473354
$ws.Range("M25").Value = -41.052631578947

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------------
Set-TextCell $ws "D26" "0" "C14"
Set-TextCell $ws "E26" "***.*" "E14"
$ws.Range("L26").Value = -20

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 4
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 55
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 77.419354838709

# ---------------------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -5
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -70.769230769230

# ---------------------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------------------
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 17
$ws.Range("K29").Value = 21.428571428571
$ws.Range("L29").Value = -10.526315789473
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -69.090909090909
